$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9136176705360413
$ws.Range("B1").Value = 1.564281940460205
$ws.Range("C1").Value = 6.000051021575928
$ws.Range("D1").Value = 1.801077961921692
$ws.Range("E1").Value = 1.093969583511353
